# Adds the 5 sample test codes used to populate the bulk-upload template's
# "Sample Code" column (rows 2-6), matching the new shared strings added in
# the target diff: VL042284501, VL042284502, VL042284512, VL042284513,
# VL0822020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sampleCodes = @(
    "VL042284501",
    "VL042284502",
    "VL042284512",
    "VL042284513",
    "VL0822020"
)

for ($i = 0; $i -lt $sampleCodes.Length; $i++) {
    $rowIndex = $i + 2
    $ws.Cells.Item($rowIndex, 1).Value = $sampleCodes[$i]
}
